# Finance workbook update: A2 becomes the numeric row-id "1" (kept in a
# text-formatted cell, matching the style already used for header cell C1),
# and two new rows of ledger data are appended (rows 3 and 4).
#
# NOTE: columns A and C of the new rows hold digit-only strings ("3", "69",
# "4", "123"). Typing them directly would be auto-coerced into numbers by
# this engine, so each one is produced as a text formula ("="3"") and then
# frozen back into a literal value via copy / paste-special-values, which
# keeps the cell's real data type as text without touching its style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Row 2: A2 was the text "X1"; it is replaced by the number 1 stored in a
# text-formatted ("@") cell, so it keeps rendering like the other text cells.
$ws.Range("A2").Value = 1
$ws.Range("A2").NumberFormat = "@"
# B2 / C2 already hold "19.01.2023" / "500" and stay untouched.

# Row 3: new ledger entry
Set-TextValue $ws.Range("A3") "3"
$ws.Range("B3").Value = "20.01.2023"
Set-TextValue $ws.Range("C3") "69"

# Row 4: new ledger entry
Set-TextValue $ws.Range("A4") "4"
$ws.Range("B4").Value = "21.01.2023"
Set-TextValue $ws.Range("C4") "123"

$ws.Range("A2").Select() | Out-Null

$wb.Application.CutCopyMode = $false
